# Add simple address handling
# ---------------------------------------------------------------------------
# This script reproduces (as closely as the COM surface allows) the diff that:
#   - inserts 6 new "studyDesign" metadata rows (+ 1 blank spacer row) above
#     the existing soa/epoch preview table on the "studyDesign" sheet
#   - widens studyDesign column A and bumps its zoom to 140%
#   - moves the active tab / selection from "soa" back to "study" (cell D2)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsStudy  = $wb.Worksheets.Item("study")
$wsDesign = $wb.Worksheets.Item("studyDesign")

# --- studyDesign: make room for the 6 new header rows ----------------------
# Old row 1 -> new row 8 (row 7 is left empty as a spacer, matching the diff).
$wsDesign.Range("A1:A7").EntireRow.Insert() | Out-Null

# --- studyDesign: column A is now wider to fit the new labels --------------
$wsDesign.Columns.Item(1).ColumnWidth = 25.1640625

# --- studyDesign: new metadata rows -----------------------------------------
# NOTE: this runtime's function-call argument binder is unreliable with
# PowerShell's named-parameter syntax (-Name value); use positional args.
function Set-LabelValueRow {
    param([string]$Row, [string]$Label, [string]$Value, [string]$ValueHAlign, [bool]$ValueTopAlign, [bool]$Italic)

    $labelCell = $wsDesign.Range("A$Row")
    $labelCell.Value = $Label
    $labelCell.Font.Bold = $true
    $labelCell.HorizontalAlignment = -4152   # xlRight
    $labelCell.VerticalAlignment = -4160     # xlTop

    $valueRange = $wsDesign.Range("B${Row}:E${Row}")
    $wsDesign.Range("B$Row").Value = $Value
    if ($Italic) {
        $valueRange.Font.Italic = $true
    }
    if ($ValueHAlign -eq "left") {
        $valueRange.HorizontalAlignment = -4131   # xlLeft
    } elseif ($ValueHAlign -eq "right") {
        $valueRange.HorizontalAlignment = -4152   # xlRight
    }
    if ($ValueTopAlign) {
        $valueRange.VerticalAlignment = -4160     # xlTop
    }
    $valueRange.Merge() | Out-Null
}

Set-LabelValueRow "1" "therapeuticAreas" "Not supported yet" "left" $false $true
Set-LabelValueRow "2" "studyDesignRationale" '"Study design rationale put here"' "left" $false $false
Set-LabelValueRow "3" "studyDesignBlindingScheme" "C49659=OPEN LABEL" "left" $false $false
Set-LabelValueRow "4" "trialIntentTypes" "C15714=BASIC SCIENCE, C139174=DEVICE FEASIBILITY" "left" $true $false
Set-LabelValueRow "5" "trialTypes" "C12345=Observational" "left" $true $false
Set-LabelValueRow "6" "interventionModel" "C12346=None" "left" $true $false

# --- studyDesign: the shifted-down "Active"/"Placebo" rows get right-aligned labels
$wsDesign.Range("A9").HorizontalAlignment = -4152   # xlRight
$wsDesign.Range("A10").HorizontalAlignment = -4152  # xlRight

# --- studyDesign: activate, zoom to 140%, select B12 ------------------------
$wsDesign.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 140
$wsDesign.Range("B12").Select() | Out-Null

# --- soa: no longer the selected tab (handled implicitly by activating
#     "study" last below), selection/zoom otherwise untouched --------------

# --- study: becomes the selected tab, with D2 as the active cell -----------
$wsStudy.Activate() | Out-Null
$wsStudy.Range("D2").Select() | Out-Null

Write-Host "Edit complete."
